# Updated cryptos list with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.290.82"
$ws.Range("E2").Value = "  -1.74%  "
$ws.Range("D3").Value = "3.751.42"
$ws.Range("E3").Value = "  +0.04%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "593.36"
$ws.Range("E5").Value = "  -0.39%  "
$ws.Range("D6").Value = "165.42"
$ws.Range("E6").Value = "  -1.23%  "
$ws.Range("D7").Value = "3.750.38"
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("E9").Value = "  -0.47%  "
$ws.Range("E10").Value = "  -0.74%  "
$ws.Range("D11").Value = "6.34"
$ws.Range("E11").Value = "  -2.08%  "
$ws.Range("D12").Value = "0.447"
$ws.Range("E12").Value = "  -0.06%  "
$ws.Range("D13").Value = "0.0000254"
$ws.Range("E13").Value = "  -2.04%  "
$ws.Range("D14").Value = "36.04"
$ws.Range("E14").Value = "  -0.05%  "
$ws.Range("D15").Value = "4.386.33"
$ws.Range("E15").Value = "  +0.15%  "
$ws.Range("D16").Value = "3.753.22"
$ws.Range("E16").Value = "  -0.13%  "
$ws.Range("D17").Value = "18.39"
$ws.Range("E17").Value = "  +2.44%  "
$ws.Range("D18").Value = "67.341.74"
$ws.Range("E18").Value = "  -1.62%  "
$ws.Range("E19").Value = "  +0.06%  "
$ws.Range("D20").Value = "6.96"
$ws.Range("E20").Value = "  -0.76%  "
$ws.Range("D21").Value = "9.95"
$ws.Range("E21").Value = "  -7.54%  "
$ws.Range("D22").Value = "453.96"
$ws.Range("E22").Value = "  -2.48%  "
$ws.Range("D23").Value = "0.694"
$ws.Range("E23").Value = "  -0.55%  "
$ws.Range("D24").Value = "0.0000153"
$ws.Range("E24").Value = "  +5.48%  "
$ws.Range("D25").Value = "83.11"
$ws.Range("E25").Value = "  -1.77%  "
$ws.Range("D26").Value = "11.85"
$ws.Range("E26").Value = "  -1.14%  "
$ws.Range("D27").Value = "2.13"
$ws.Range("E27").Value = "  -2.62%  "
$ws.Range("D28").Value = "10.09"
$ws.Range("E28").Value = "  +0.88%  "
$ws.Range("E29").Value = "  +0.03%  "
$ws.Range("D30").Value = "2.77"
$ws.Range("E30").Value = "  -0.45%  "
$ws.Range("D31").Value = "7.24"
$ws.Range("E31").Value = "  -1.00%  "
$ws.Range("D32").Value = "29.56"
$ws.Range("E32").Value = "  -0.96%  "
$ws.Range("D33").Value = "2.17"
$ws.Range("E33").Value = "  -0.30%  "
$ws.Range("D34").Value = "9.16"
$ws.Range("E34").Value = "  -0.59%  "
$ws.Range("E35").Value = "  -0.19%  "
$ws.Range("D36").Value = "3.706.95"
$ws.Range("E36").Value = "  +0.05%  "
$ws.Range("D37").Value = "0.0998"
$ws.Range("E37").Value = "  -1.15%  "
$ws.Range("D38").Value = "3.31"
$ws.Range("E38").Value = "  -1.90%  "
$ws.Range("D39").Value = "0.137"
$ws.Range("E39").Value = "  -1.15%  "
$ws.Range("D40").Value = "0.993"
$ws.Range("E40").Value = "  -0.43%  "
$ws.Range("D41").Value = "5.73"
$ws.Range("E41").Value = "  -1.54%  "
$ws.Range("E42").Value = "  +0.20%  "
$ws.Range("D44").Value = "45.18"
$ws.Range("E44").Value = "  +2.90%  "
$ws.Range("D45").Value = "0.297"
$ws.Range("E45").Value = "  -2.30%  "
$ws.Range("D46").Value = "46.94"
$ws.Range("E46").Value = "  +2.14%  "

# Rows 47 and 48 swap: Cosmos <-> Monero (with their own updated price/volume)
$ws.Range("B47").Value = "Monero"
$ws.Range("C47").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D47").Value = "148.55"
$ws.Range("E47").Value = "  +1.14%  "

$ws.Range("B48").Value = "Cosmos"
$ws.Range("C48").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.30"
$ws.Range("E48").Value = "  -3.06%  "

$ws.Range("D49").Value = "1.83"
$ws.Range("E49").Value = "  -4.90%  "
$ws.Range("D50").Value = "388.54"
$ws.Range("E50").Value = "  -0.35%  "
$ws.Range("D51").Value = "25.95"
$ws.Range("E51").Value = "  +0.84%  "
